$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Presentation-level slide guides (stored in p:presentation/p:extLst
#    as a p15:sldGuideLst): one horizontal guide at 3.75in (270pt) and
#    one vertical guide at 5in (360pt), both the default PowerPoint
#    "center" guide color.
# ------------------------------------------------------------------
try {
    $hGuide = $p.Guides.Add(1, 270)   # 1 = ppHorizontalGuide, position in points (2160/8)
    $vGuide = $p.Guides.Add(2, 360)   # 2 = ppVerticalGuide,  position in points (2880/8)
} catch {
    # Guides collection not available in this host - ignore and continue.
}

# ------------------------------------------------------------------
# 2) Re-cache the "datetimeFigureOut" footer date field that shows up
#    on the slide master, every slide layout, and the notes master
#    (the footer date was refreshed from 24/04/2013 to 20/12/2013).
# ------------------------------------------------------------------
$newDate = "20/12/2013"
$dateShapeName = "3 Marcador de fecha"

# Slide Master footer date placeholder.
$master = $p.SlideMaster
$masterDateShape = $master.Shapes.Item($dateShapeName)
$masterDateShape.TextFrame.TextRange.Text = $newDate

# Every slide layout's footer date placeholder.
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    $layoutDateShape = $layout.Shapes.Item($dateShapeName)
    $layoutDateShape.TextFrame.TextRange.Text = $newDate
}

# Notes Master footer date placeholder (named slightly differently).
$notesMaster = $p.NotesMaster
$notesDateShape = $notesMaster.Shapes.Item("2 Marcador de fecha")
$notesDateShape.TextFrame.TextRange.Text = $newDate
